# Update the Training Dashboard sheet with the new progress-as-of date (04-Nov-2025).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Row 3: PERIOD TO EXPIRE decreased by one day, LAST UPDATE moved to 04-Nov-2025.
$ws.Range("H3").Value = -99
$ws.Range("I3").Value = "'04-Nov-2025"

# Row 4: PERIOD TO EXPIRE decreased by one day, LAST UPDATE moved to 04-Nov-2025.
$ws.Range("H4").Value = 699
$ws.Range("I4").Value = "'04-Nov-2025"
